$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "ValidLogin" and add a new sheet "InvalidLogin" after it.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

# Valid-login credentials (ValidLogin sheet) - enter the data row first, then the headers,
# matching the order the values were typed in.
$ws1.Range("A2").Value = "Admin"
$ws1.Range("B2").Value = "admin123"
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"

# Invalid-login credentials (InvalidLogin sheet)
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# Leave ValidLogin's used range selected, and make InvalidLogin the active sheet with B3 selected.
[void]$ws1.Range("A1:B2").Select()
[void]$ws2.Range("B3").Select()
[void]$ws2.Activate()
